$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_study")

# Study Person Last Name / First Name values
$ws.Range("B50").Value = "Krumsieck"
$ws.Range("B51").Value = "Jens"

# Study Person Affiliation value
$ws.Range("B57").Value = "Johann Heinrich von Thünen-Institut, Zentrum für Informationsmanagement;Technische Universität Braunschweig, Institut für Anorganische und Analytische Chemie"

# New row 61: Comment[ORCID]
$ws.Range("A61").Value = "Comment[ORCID]"
$ws.Range("B61").Value = "0000-0001-6242-5846"
